$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''318.12'
$ws.Range("E2").Value = '''4.28%'
$ws.Range("G2").Value = '''15'
$ws.Range("D2:G2").Style = "Normal"

$ws.Range("D3").Value = '''39.48'
$ws.Range("E3").Value = '''2.81%'
$ws.Range("G3").Value = '''15'
$ws.Range("D3:G3").Style = "Normal"

$ws.Range("D4").Value = '''5.126'
$ws.Range("E4").Value = '''0.69%'
$ws.Range("G4").Value = '''15'
$ws.Range("D4:G4").Style = "Normal"

$ws.Range("D5").Value = '''0.08204'
$ws.Range("E5").Value = '''1.78%'
$ws.Range("G5").Value = '''15'
$ws.Range("D5:G5").Style = "Normal"

$ws.Range("E6").Value = '''5.24%'
$ws.Range("G6").Value = '''15'
$ws.Range("E6:G6").Style = "Normal"

$ws.Range("D7").Value = '''8.282'
$ws.Range("E7").Value = '''4.19%'
$ws.Range("G7").Value = '''15'
$ws.Range("D7:G7").Style = "Normal"

$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = '''4.285'
$ws.Range("E8").Value = '''2.24%'
$ws.Range("G8").Value = '''15'
$ws.Range("B8:G8").Style = "Normal"

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9328'
$ws.Range("E9").Value = '''0.31%'
$ws.Range("G9").Value = '''15'
$ws.Range("B9:G9").Style = "Normal"

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1406'
$ws.Range("E10").Value = '''-3.86%'
$ws.Range("G10").Value = '''15'
$ws.Range("B10:G10").Style = "Normal"

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1992'
$ws.Range("E11").Value = '''3.65%'
$ws.Range("G11").Value = '''15'
$ws.Range("B11:G11").Style = "Normal"

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09144'
$ws.Range("E12").Value = '''1.62%'
$ws.Range("G12").Value = '''15'
$ws.Range("B12:G12").Style = "Normal"

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03558'
$ws.Range("E13").Value = '''0.98%'
$ws.Range("G13").Value = '''15'
$ws.Range("B13:G13").Style = "Normal"

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09826'
$ws.Range("E14").Value = '''0.44%'
$ws.Range("G14").Value = '''15'
$ws.Range("B14:G14").Style = "Normal"

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001402'
$ws.Range("E15").Value = '''-0.07%'
$ws.Range("G15").Value = '''15'
$ws.Range("B15:G15").Style = "Normal"

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006382'
$ws.Range("E16").Value = '''7.67%'
$ws.Range("G16").Value = '''15'
$ws.Range("B16:G16").Style = "Normal"

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.660'
$ws.Range("E17").Value = '''-1.59%'
$ws.Range("G17").Value = '''15'
$ws.Range("B17:G17").Style = "Normal"

$ws.Range("D18").Value = '''3.279'
$ws.Range("E18").Value = '''-4.15%'
$ws.Range("G18").Value = '''15'
$ws.Range("D18:G18").Style = "Normal"

$ws.Range("D19").Value = '''0.3461'
$ws.Range("E19").Value = '''-0.02%'
$ws.Range("G19").Value = '''15'
$ws.Range("D19:G19").Style = "Normal"

$ws.Range("D20").Value = '''0.1304'
$ws.Range("E20").Value = '''-0.67%'
$ws.Range("G20").Value = '''15'
$ws.Range("D20:G20").Style = "Normal"

$ws.Range("D21").Value = '''4.906'
$ws.Range("E21").Value = '''1.72%'
$ws.Range("G21").Value = '''15'
$ws.Range("D21:G21").Style = "Normal"

$ws.Range("D22").Value = '''0.2451'
$ws.Range("G22").Value = '''15'
$ws.Range("D22:G22").Style = "Normal"

$ws.Range("D23").Value = '''0.04327'
$ws.Range("E23").Value = '''-0.94%'
$ws.Range("G23").Value = '''15'
$ws.Range("D23:G23").Style = "Normal"

$ws.Range("D24").Value = '''0.001225'
$ws.Range("E24").Value = '''-0.67%'
$ws.Range("G24").Value = '''15'
$ws.Range("D24:G24").Style = "Normal"

$ws.Range("D25").Value = '''0.004786'
$ws.Range("E25").Value = '''12.61%'
$ws.Range("G25").Value = '''15'
$ws.Range("D25:G25").Style = "Normal"

$ws.Range("D26").Value = '''0.0001300'
$ws.Range("E26").Value = '''-0.15%'
$ws.Range("G26").Value = '''15'
$ws.Range("D26:G26").Style = "Normal"

$ws.Range("D27").Value = '''0.0004000'
$ws.Range("E27").Value = '''-10.07%'
$ws.Range("G27").Value = '''15'
$ws.Range("D27:G27").Style = "Normal"

$ws.Range("G28").Value = '''15'
$ws.Range("G28").Style = "Normal"

$ws.Range("G29").Value = '''15'
$ws.Range("G29").Style = "Normal"

$ws.Range("G30").Value = '''15'
$ws.Range("G30").Style = "Normal"

$ws.Range("G31").Value = '''15'
$ws.Range("G31").Style = "Normal"

$ws.Range("G32").Value = '''15'
$ws.Range("G32").Style = "Normal"

$ws.Range("G33").Value = '''15'
$ws.Range("G33").Style = "Normal"

$ws.Range("G34").Value = '''15'
$ws.Range("G34").Style = "Normal"

$ws.Range("G35").Value = '''15'
$ws.Range("G35").Style = "Normal"

$ws.Range("G36").Value = '''15'
$ws.Range("G36").Style = "Normal"

$ws.Range("G37").Value = '''15'
$ws.Range("G37").Style = "Normal"

$ws.Range("G38").Value = '''15'
$ws.Range("G38").Style = "Normal"

$ws.Range("D39").Value = '''0.02229'
$ws.Range("E39").Value = '''8.22%'
$ws.Range("G39").Value = '''15'
$ws.Range("D39:G39").Style = "Normal"

$ws.Range("D40").Value = '''0.05266'
$ws.Range("E40").Value = '''4.83%'
$ws.Range("G40").Value = '''15'
$ws.Range("D40:G40").Style = "Normal"

$ws.Range("D41").Value = '''0.007529'
$ws.Range("E41").Value = '''0.73%'
$ws.Range("G41").Value = '''15'
$ws.Range("D41:G41").Style = "Normal"

$ws.Range("D42").Value = '''0.009763'
$ws.Range("E42").Value = '''-3.54%'
$ws.Range("G42").Value = '''15'
$ws.Range("D42:G42").Style = "Normal"

$ws.Range("D43").Value = '''0.1382'
$ws.Range("E43").Value = '''2.64%'
$ws.Range("G43").Value = '''15'
$ws.Range("D43:G43").Style = "Normal"

$ws.Range("D44").Value = '''0.002150'
$ws.Range("E44").Value = '''0.32%'
$ws.Range("G44").Value = '''15'
$ws.Range("D44:G44").Style = "Normal"

$ws.Range("D45").Value = '''0.009865'
$ws.Range("E45").Value = '''10.25%'
$ws.Range("G45").Value = '''15'
$ws.Range("D45:G45").Style = "Normal"

$ws.Range("D46").Value = '''0.00006452'
$ws.Range("E46").Value = '''4.10%'
$ws.Range("G46").Value = '''15'
$ws.Range("D46:G46").Style = "Normal"

$ws.Range("E47").Value = '''-0.14%'
$ws.Range("G47").Value = '''15'
$ws.Range("E47:G47").Style = "Normal"

$ws.Range("D48").Value = '''0.002769'
$ws.Range("E48").Value = '''-1.86%'
$ws.Range("G48").Value = '''15'
$ws.Range("D48:G48").Style = "Normal"

$ws.Range("E49").Value = '''-25.12%'
$ws.Range("G49").Value = '''15'
$ws.Range("E49:G49").Style = "Normal"

$ws.Range("D50").Value = '''0.00002100'
$ws.Range("E50").Value = '''-0.14%'
$ws.Range("G50").Value = '''15'
$ws.Range("D50:G50").Style = "Normal"

$ws.Range("D51").Value = '''0.0002000'
$ws.Range("E51").Value = '''-0.14%'
$ws.Range("G51").Value = '''15'
$ws.Range("D51:G51").Style = "Normal"
